$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh generation Date/time ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-05T11:54:16+00:00"

# --- Elements sheet: AuthorDocumentEntry.person -> AuthorDocumentEntry.person[x] row ---
$ws = $wb.Worksheets.Item("Elements")

# ID / Path columns (A4, B4)
$ws.Range("A4").Value = "AuthorDocumentEntry.person[x]"
$ws.Range("B4").Value = "AuthorDocumentEntry.person[x]"

# Type(s) column (K4): drop the "Reference(...|...)" wrapper/pipes, join with newline then concat
$newline = [char]10
$typesText = "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorPS" + $newline + "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorPatienthttps://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorSNRhttps://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorSystem"
$ws.Range("K4").Value = $typesText

# Definition column (M4): "author" -> "Author"
$ws.Range("M4").Value = "Author"

# Base Path column (AF4): "Author.person" -> "Author.person[x]"
$ws.Range("AF4").Value = "Author.person[x]"

# Re-fit the Type(s) column now that its longest entry is shorter than before
$ws.Columns.Item(11).AutoFit()
